# Apply updated dSF (column F) values pulled from a repull of the data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -8
    3  = -7
    5  = -7
    6  = -3
    8  = 9
    10 = 9
    12 = 7
    13 = 0
    17 = 6
    20 = 1
    21 = 1
    23 = -2
    25 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
